# Generate Report for Handback
# - Update "Ready for handoff" status text to "Handed back: in sync with en-US"
#   everywhere it is used (Overview + per-language sheets pick it up via the
#   shared string, so we just rewrite the cells that hold that text).
# - Populate the "Latest Target File" (F) and "Latest Handback File" (G)
#   columns on the zh-cn / de-de sheets with hyperlinked file names, mirroring
#   the existing "Source File Name" (A) / "Latest Handoff File" (D) links.
# - Stamp the "Latest Handback DateTime" (H) column with the real handback
#   timestamp now that the handback has happened.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

# ---- Overview sheet: refresh the per-language status columns ----
$ov = $wb.Worksheets.Item("Overview")
for ($r = 2; $r -le 3; $r++) {
    if ($ov.Cells.Item($r, 2).Value2 -eq $oldStatus) {
        $ov.Cells.Item($r, 2).Value = $newStatus
    }
    if ($ov.Cells.Item($r, 3).Value2 -eq $oldStatus) {
        $ov.Cells.Item($r, 3).Value = $newStatus
    }
}

# ---- zh-cn sheet ----
$zh = $wb.Worksheets.Item("zh-cn")

# Status column (C) reflects the handback too.
if ($zh.Range("C2").Value2 -eq $oldStatus) { $zh.Range("C2").Value = $newStatus }
if ($zh.Range("C3").Value2 -eq $oldStatus) { $zh.Range("C3").Value = $newStatus }

# Row 2 (882fc755 file): populate Latest Target File (F) / Latest Handback File (G)
$zh.Range("F2").Value = "882fc755-095b-4cc6-af4b-658cf8c09ce0.md"
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/8703c7e968ff5d6df87973249e5a5e12abce2ec8/e2e/882fc755-095b-4cc6-af4b-658cf8c09ce0.md", "", "", "882fc755-095b-4cc6-af4b-658cf8c09ce0.md") | Out-Null

$zh.Range("G2").Value = "882fc755-095b-4cc6-af4b-658cf8c09ce0.c33c0f7ceb83615306440c6713b62ebe27099bc9.zh-cn.xlf"
$zh.Hyperlinks.Add($zh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/15636e0467fa27a3a42894f0ba9e7133efec702b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/882fc755-095b-4cc6-af4b-658cf8c09ce0.c33c0f7ceb83615306440c6713b62ebe27099bc9.zh-cn.xlf", "", "", "882fc755-095b-4cc6-af4b-658cf8c09ce0.c33c0f7ceb83615306440c6713b62ebe27099bc9.zh-cn.xlf") | Out-Null

# Row 3 (9b8ec366 file): populate Latest Target File (F) / Latest Handback File (G)
$zh.Range("F3").Value = "9b8ec366-ef01-4bbb-b031-849c0146b210.md"
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/8703c7e968ff5d6df87973249e5a5e12abce2ec8/e2e/9b8ec366-ef01-4bbb-b031-849c0146b210.md", "", "", "9b8ec366-ef01-4bbb-b031-849c0146b210.md") | Out-Null

$zh.Range("G3").Value = "9b8ec366-ef01-4bbb-b031-849c0146b210.0516ba58e908ea33dab883b9719b8860793ac75e.zh-cn.xlf"
$zh.Hyperlinks.Add($zh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/15636e0467fa27a3a42894f0ba9e7133efec702b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/9b8ec366-ef01-4bbb-b031-849c0146b210.0516ba58e908ea33dab883b9719b8860793ac75e.zh-cn.xlf", "", "", "9b8ec366-ef01-4bbb-b031-849c0146b210.0516ba58e908ea33dab883b9719b8860793ac75e.zh-cn.xlf") | Out-Null

# zh-cn handback finished at 2016-03-24 01:09:05
$zh.Range("H2").Value = "2016-03-24 01:09:05"
$zh.Range("H3").Value = "2016-03-24 01:09:05"

# ---- de-de sheet ----
$de = $wb.Worksheets.Item("de-de")

if ($de.Range("C2").Value2 -eq $oldStatus) { $de.Range("C2").Value = $newStatus }
if ($de.Range("C3").Value2 -eq $oldStatus) { $de.Range("C3").Value = $newStatus }

# Row 2 (882fc755 file): populate Latest Target File (F) / Latest Handback File (G)
$de.Range("F2").Value = "882fc755-095b-4cc6-af4b-658cf8c09ce0.md"
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/8703c7e968ff5d6df87973249e5a5e12abce2ec8/e2e/882fc755-095b-4cc6-af4b-658cf8c09ce0.md", "", "", "882fc755-095b-4cc6-af4b-658cf8c09ce0.md") | Out-Null

$de.Range("G2").Value = "882fc755-095b-4cc6-af4b-658cf8c09ce0.c33c0f7ceb83615306440c6713b62ebe27099bc9.de-de.xlf"
$de.Hyperlinks.Add($de.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9b644cb9e11c21af7c82b206271f079f0f50aa31/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/882fc755-095b-4cc6-af4b-658cf8c09ce0.c33c0f7ceb83615306440c6713b62ebe27099bc9.de-de.xlf", "", "", "882fc755-095b-4cc6-af4b-658cf8c09ce0.c33c0f7ceb83615306440c6713b62ebe27099bc9.de-de.xlf") | Out-Null

# Row 3 (9b8ec366 file): populate Latest Target File (F) / Latest Handback File (G)
$de.Range("F3").Value = "9b8ec366-ef01-4bbb-b031-849c0146b210.md"
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/8703c7e968ff5d6df87973249e5a5e12abce2ec8/e2e/9b8ec366-ef01-4bbb-b031-849c0146b210.md", "", "", "9b8ec366-ef01-4bbb-b031-849c0146b210.md") | Out-Null

$de.Range("G3").Value = "9b8ec366-ef01-4bbb-b031-849c0146b210.0516ba58e908ea33dab883b9719b8860793ac75e.de-de.xlf"
$de.Hyperlinks.Add($de.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9b644cb9e11c21af7c82b206271f079f0f50aa31/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/9b8ec366-ef01-4bbb-b031-849c0146b210.0516ba58e908ea33dab883b9719b8860793ac75e.de-de.xlf", "", "", "9b8ec366-ef01-4bbb-b031-849c0146b210.0516ba58e908ea33dab883b9719b8860793ac75e.de-de.xlf") | Out-Null

# de-de handback finished at 2016-03-24 01:09:12 (distinct from zh-cn's time)
$de.Range("H2").Value = "2016-03-24 01:09:12"
$de.Range("H3").Value = "2016-03-24 01:09:12"

Write-Host "Handback report generated."
